$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row update
$ws.Range("C1").Value = "prediction"
$ws.Range("D1").Value = "rejection-f"
$ws.Range("E1").Value = "max"

# Data row update
$ws.Range("C2").Value = "s__Turicimonas sp900542195"
$ws.Range("D2").Value = "s__Turicimonas sp900542195"
$ws.Range("E2").Value = 1
